$d = $word.ActiveDocument

# First paragraph: the hidden "**ID__...__ID**" marker paragraph.
$p1 = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right, each w:space="5") and
# change the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

$p1.Range.ParagraphFormat.LeftIndent = 11.25

# Replace the marker text (and absorb the trailing space run into the
# same run) with the updated subpart-scoped marker id.
$rng = $p1.Range
$rng.Find.Execute("**ID__AFFARS_5323_topic_9__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_SUBPART_5323_90__ID**", 2)

Write-Host "done"
